$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Vendor info block updates (rows 3-7) ---

# Vendor Name / Address (row 3) - new vendor
$ws.Range("B3").Value = "2GO Express"
$ws.Range("D3").Value = "BREDCO, Port 2, Reclamation Area, Brgy. 10, Bacolod City"

# Phone Number (row 4) - new value; Fax Number value cleared (fax number removed)
$ws.Range("B4").Value = "(034) 704-1339"
$ws.Range("D4").ClearContents()

# Row 5 was Fax Number(blank)/Email - now becomes Email:/Contact Person:
$ws.Range("A5").Value = "Email:"
$ws.Range("C5").Value = "Contact Person:"
$ws.Range("D5").ClearContents()

# Row 6 was Contact Person/Terms - now becomes Terms:/Type:
$ws.Range("A6").Value = "Terms:"
$ws.Range("C6").Value = "Type:"
$ws.Range("D6").ClearContents()

# Row 7 was Notes: (merged B7:D7, blank) - now becomes EWT(%): / 0 / Notes:
$ws.Range("B7:D7").UnMerge()
$ws.Range("A7").Value = "EWT(%):"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "Notes:"
$ws.Range("D7").ClearContents()

# --- Remove the item list rows (11-22), keep header rows 9-10 ---
$ws.Range("A11:D22").EntireRow.Delete()

# Update selection to the new last cell
[void]$ws.Range("D10").Select()
